# "Add files via upload" — refreshed export of Saldo.xlsx.
#
# The new export has one additional data row (account 008007764 / LUIS /
# 10115.23) and refreshed balances/names for most of the other rows. The
# trailing "Filtros aplicados..." note slides from row 12 to row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 7 ("PEDRO") so everything from
# there down (including the filter-notes row) shifts down by one, matching
# the new row count (A1:C13 instead of A1:C12).
$ws.Rows.Item(7).Insert()

# Account numbers are zero-padded strings ("008030888", ...). Mark column A
# as Text before writing them so Excel keeps the leading zeros instead of
# reinterpreting them as numbers.
$ws.Range("A2:A11").NumberFormat = "@"

# Final table contents (account, name, balance) for rows 2-11.
$data = @(
    @("008030888", "SONIA",      400000),
    @("005547467", "CARLA",      171875),
    @("004546050", "LUIS",       20307.5),
    @("005142661", "SABRINA",    16170.39),
    @("004581652", "CINCO",      14455.12),
    @("008007764", "LUIS",       10115.23),
    @("004214592", "MERG",       5727.04),
    @("004240400", "ADRIANA",    3000),
    @("001761119", "BLUEMETRIX", 1734.91),
    @("004504449", "KELMA",      1000)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Extend the (hidden) filter-database range to include the new row.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Export!`$A`$1:`$C`$11"

# Reposition the window / selection the way the saved file shows it.
$win = $wb.Windows.Item(1)
$win.Left = 28815
$win.Top = -16320
$ws.Range("B18").Select()
